# Scheduled-runner refresh of the per-profession "currentAveragePrice*" /
# "LevePrice*" / "LeveProfit*" market-data columns (H:N) across all eight
# job sheets. Values are literal (no formulas in this workbook), so each
# touched cell is simply re-written with the freshly captured figure; a
# handful of cells whose leves no longer price out (or newly do) are
# cleared/added to match.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 335
$ws.Range("I2").Value = 533.3333
$ws.Range("J2").Value = 136.66667
$ws.Range("K2").Value = 533.3333
$ws.Range("L2").Value = 136.66667
$ws.Range("M2").Value = -420.3333
$ws.Range("N2").Value = -362.66667

$ws.Range("H5").Value = 508
$ws.Range("I5").Value = 460
$ws.Range("J5").Value = 700
$ws.Range("K5").Value = 460
$ws.Range("L5").Value = 700
$ws.Range("M5").Value = -345
$ws.Range("N5").Value = -930

$ws.Range("H18").Value = 10608.723
$ws.Range("I18").Value = 2963.3333
$ws.Range("J18").Value = 14431.417
$ws.Range("K18").Value = 2963.3333
$ws.Range("L18").Value = 14431.417
$ws.Range("M18").Value = -2679.3333
$ws.Range("N18").Value = -14999.417

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null

$ws.Range("H98").Value = 1227.4
$ws.Range("I98").Value = 1258.409
$ws.Range("K98").Value = 1258.409
$ws.Range("M98").Value = 239.5909999999999

$ws.Range("H112").Value = 6255.222
$ws.Range("I112").Value = 398
$ws.Range("K112").Value = 1194
$ws.Range("M112").Value = -86

$ws.Range("H122").Value = 1227.4
$ws.Range("I122").Value = 1258.409
$ws.Range("K122").Value = 3775.227
$ws.Range("M122").Value = -1325.227

$ws.Range("H141").Value = 2337847.2
$ws.Range("I141").Value = 4670013.5
$ws.Range("K141").Value = 14010040.5
$ws.Range("M141").Value = -14004860.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3436.07
$ws.Range("I32").Value = 2061.6326
$ws.Range("K32").Value = 2061.6326
$ws.Range("M32").Value = -1774.6326

$ws.Range("H122").Value = 3201
$ws.Range("I122").Value = 3401.375
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 10204.125
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("M122").Value = -7754.125
$ws.Range("N122").Value = -12900.0001

$ws.Range("H132").Value = 1306.3469
$ws.Range("J132").Value = 1741.3478
$ws.Range("L132").Value = 5224.0434
$ws.Range("N132").Value = -10284.0434

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 83068.89
$ws.Range("I86").Value = 1530.1177
$ws.Range("J86").Value = 221684.8
$ws.Range("K86").Value = 1530.1177
$ws.Range("L86").Value = 221684.8
$ws.Range("M86").Value = -407.1177
$ws.Range("N86").Value = -223930.8

$ws.Range("H89").Value = 83068.89
$ws.Range("I89").Value = 1530.1177
$ws.Range("J89").Value = 221684.8
$ws.Range("K89").Value = 7650.5885
$ws.Range("L89").Value = 1108424
$ws.Range("M89").Value = -2034.5885
$ws.Range("N89").Value = -1119656

$ws.Range("H94").Value = 613.52
$ws.Range("I94").Value = 583.86365
$ws.Range("K94").Value = 583.86365
$ws.Range("M94").Value = -132.86365

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 10000
$ws.Range("J17").Value = 10000
$ws.Range("L17").Value = 10000
$ws.Range("N17").Value = -10348

$ws.Range("H31").Value = 2492.9792
$ws.Range("I31").Value = 1663.1
$ws.Range("K31").Value = 1663.1
$ws.Range("M31").Value = -1368.1

$ws.Range("H34").Value = 2492.9792
$ws.Range("I34").Value = 1663.1
$ws.Range("K34").Value = 1663.1
$ws.Range("M34").Value = -1461.1

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = $null

$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = $null

$ws.Range("H60").Value = 11414
$ws.Range("J60").Value = 11414
$ws.Range("L60").Value = 11414
$ws.Range("N60").Value = -12436

$ws.Range("H68").Value = 24700

$ws.Range("H71").Value = 24700

$ws.Range("H74").Value = 31000
$ws.Range("J74").Value = 31000
$ws.Range("L74").Value = 31000
$ws.Range("N74").Value = -32748

$ws.Range("H77").Value = 31000
$ws.Range("J77").Value = 31000
$ws.Range("L77").Value = 93000
$ws.Range("N77").Value = -101736

$ws.Range("H99").Value = 669100.5600000001
$ws.Range("I99").Value = 1668283.1
$ws.Range("J99").Value = 2978.7778
$ws.Range("K99").Value = 1668283.1
$ws.Range("L99").Value = 2978.7778
$ws.Range("M99").Value = -1666785.1
$ws.Range("N99").Value = -5974.7778

$ws.Range("H126").Value = 669100.5600000001
$ws.Range("I126").Value = 1668283.1
$ws.Range("J126").Value = 2978.7778
$ws.Range("K126").Value = 5004849.300000001
$ws.Range("L126").Value = 8936.3334
$ws.Range("M126").Value = -5002379.300000001
$ws.Range("N126").Value = -13876.3334

$ws.Range("H132").Value = 1907.0857
$ws.Range("I132").Value = 1198.2
$ws.Range("J132").Value = 2852.2666
$ws.Range("K132").Value = 3594.6
$ws.Range("L132").Value = 8556.799800000001
$ws.Range("M132").Value = -1064.6
$ws.Range("N132").Value = -13616.7998

$ws.Range("H134").Value = 923.41174
$ws.Range("I134").Value = 892.1667
$ws.Range("J134").Value = 998.4
$ws.Range("K134").Value = 2676.5001
$ws.Range("L134").Value = 2995.2
$ws.Range("M134").Value = -141.5001000000002
$ws.Range("N134").Value = -8065.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 167.375
$ws.Range("I2").Value = 247.25
$ws.Range("J2").Value = 87.5
$ws.Range("K2").Value = 1483.5
$ws.Range("L2").Value = 525
$ws.Range("M2").Value = -1370.5
$ws.Range("N2").Value = -751

$ws.Range("H5").Value = 743.44446
$ws.Range("J5").Value = 817.4
$ws.Range("L5").Value = 2452.2
$ws.Range("N5").Value = -2676.2

$ws.Range("H17").Value = 10050.5
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 10050.5
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 30151.5
$ws.Range("M17").Value = $null
$ws.Range("N17").Value = -30489.5

$ws.Range("H32").Value = 880
$ws.Range("I32").Value = 700
$ws.Range("J32").Value = 940
$ws.Range("K32").Value = 2100
$ws.Range("L32").Value = 2820
$ws.Range("M32").Value = -1817
$ws.Range("N32").Value = -3386

$ws.Range("H131").Value = 11061.456
$ws.Range("J131").Value = 12452.583
$ws.Range("L131").Value = 37357.749
$ws.Range("N131").Value = -47437.749

$ws.Range("H135").Value = 743.44446
$ws.Range("J135").Value = 817.4
$ws.Range("L135").Value = 7356.599999999999
$ws.Range("N135").Value = -12426.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4537
$ws.Range("I70").Value = 4401.1665
$ws.Range("K70").Value = 4401.1665
$ws.Range("M70").Value = -4131.1665

$ws.Range("H73").Value = 4537
$ws.Range("I73").Value = 4401.1665
$ws.Range("K73").Value = 4401.1665
$ws.Range("M73").Value = -3465.1665

$ws.Range("H102").Value = 3374.7693
$ws.Range("J102").Value = 2626.625
$ws.Range("L102").Value = 2626.625
$ws.Range("N102").Value = -5870.625

$ws.Range("H122").Value = 1916
$ws.Range("I122").Value = 1566.6666
$ws.Range("K122").Value = 4699.9998
$ws.Range("M122").Value = -2249.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6686.091
$ws.Range("I7").Value = 3900.75
$ws.Range("J7").Value = 8277.714
$ws.Range("K7").Value = 3900.75
$ws.Range("L7").Value = 8277.714
$ws.Range("M7").Value = -3788.75
$ws.Range("N7").Value = -8501.714

$ws.Range("H68").Value = 2620.4285
$ws.Range("I68").Value = 1960
$ws.Range("J68").Value = 3501
$ws.Range("K68").Value = 1960
$ws.Range("L68").Value = 3501
$ws.Range("M68").Value = -1211
$ws.Range("N68").Value = -4999

$ws.Range("H71").Value = 2620.4285
$ws.Range("I71").Value = 1960
$ws.Range("J71").Value = 3501
$ws.Range("K71").Value = 9800
$ws.Range("L71").Value = 17505
$ws.Range("M71").Value = -6056
$ws.Range("N71").Value = -24993

$ws.Range("H122").Value = 4909
$ws.Range("I122").Value = 1380.6
$ws.Range("K122").Value = 4141.799999999999
$ws.Range("M122").Value = -1691.799999999999

$ws.Range("H126").Value = 6686.091
$ws.Range("I126").Value = 3900.75
$ws.Range("J126").Value = 8277.714
$ws.Range("K126").Value = 11702.25
$ws.Range("L126").Value = 24833.142
$ws.Range("M126").Value = -9232.25
$ws.Range("N126").Value = -29773.142

$ws.Range("H136").Value = 2207.9644
$ws.Range("I136").Value = 1017.94446
$ws.Range("J136").Value = 4350
$ws.Range("K136").Value = 3053.83338
$ws.Range("L136").Value = 13050
$ws.Range("M136").Value = -503.83338
$ws.Range("N136").Value = -18150

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 470069.25
$ws.Range("I122").Value = 470069.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1410207.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1407757.75
$ws.Range("N122").Value = $null

$ws.Range("H136").Value = 14246929
$ws.Range("I136").Value = 19158734
$ws.Range("J136").Value = 2694.7
$ws.Range("K136").Value = 57476202
$ws.Range("L136").Value = 8084.099999999999
$ws.Range("M136").Value = -57473652
$ws.Range("N136").Value = -13184.1
